$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "OPQA-1409"

$ws.Columns.Item(1).ColumnWidth = 10.166666666666666

$ws.Range("A3").Select() | Out-Null
